$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column keeps its original text semantics (several values look
# like numbers, e.g. '1.000', '0.9999', '25.949.81') - force text format so
# Excel COM does not silently coerce them into numeric cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.949.81'
$ws.Range("E2").Value = '  -0.81%  '
$ws.Range("D3").Value = '1.745.98'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '248.72'
$ws.Range("E5").Value = '  +5.01%  '
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.5051'
$ws.Range("E7").Value = '  -6.05%  '
$ws.Range("E8").Value = '  -1.99%  '
$ws.Range("D9").Value = '0.06189'
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("D10").Value = '0.07270'
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("D11").Value = '1.741.98'
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").Value = '0.6551'
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("D13").Value = '15.19'
$ws.Range("E13").Value = '  -1.71%  '
$ws.Range("D14").Value = '4.657'
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").Value = '77.76'
$ws.Range("E15").Value = '  -0.94%  '
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '25.969.41'
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("D19").Value = '11.87'
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("E20").Value = '  +1.24%  '
$ws.Range("D21").Value = '1.968.22'
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").Value = '4.455'
$ws.Range("E22").Value = '  +2.49%  '
$ws.Range("D23").Value = '8.729'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").Value = '5.403'
$ws.Range("E24").Value = '  +2.85%  '
$ws.Range("D25").Value = '136.73'
$ws.Range("E25").Value = '  -2.02%  '
$ws.Range("D26").Value = '1.504'
$ws.Range("E26").Value = '  -0.99%  '
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("D28").Value = '1.783'
$ws.Range("E28").Value = '  -1.10%  '
$ws.Range("D29").Value = '105.70'
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("D30").Value = '3.869'
$ws.Range("E30").Value = '  +2.03%  '
$ws.Range("D31").Value = '0.08210'
$ws.Range("E31").Value = '  -1.52%  '
$ws.Range("D32").Value = '3.653'
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").Value = '0.04684'
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("E34").Value = '  +0.31%  '
$ws.Range("D35").Value = '0.9968'
$ws.Range("E35").Value = '  -1.28%  '
$ws.Range("E36").Value = '  -1.98%  '
$ws.Range("D37").Value = '2.753'
$ws.Range("E37").Value = '  +1.55%  '
$ws.Range("D38").Value = '0.01613'
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("D40").Value = '0.9997'
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").Value = '100.75'
$ws.Range("D42").Value = '0.3943'
$ws.Range("E42").Value = '  +0.23%  '
$ws.Range("D43").Value = '0.7604'
$ws.Range("E43").Value = '  +1.29%  '
$ws.Range("D44").Value = '5.011'
$ws.Range("E44").Value = '  -1.60%  '
$ws.Range("D45").Value = '0.1151'
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("D46").Value = '6.339'
$ws.Range("E46").Value = '  -0.43%  '
$ws.Range("D47").Value = '55.81'
$ws.Range("E47").Value = '  +1.91%  '
$ws.Range("D48").Value = '0.05283'
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("D49").Value = '30.70'
$ws.Range("E49").Value = '  -0.98%  '
$ws.Range("D50").Value = '7.603'
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("D51").Value = '0.3449'
$ws.Range("E51").Value = '  -1.16%  '
